# Remove the rows that were dropped from the "Export" sheet.
# Row numbers below are 1-based worksheet rows (row 1 is the header row:
# Conta | Nome | Saldo).
#
# Deleted rows (from the original workbook, before any shifting):
#   Row 4   -> 004976625 | Norton   | 80127.16
#   Row 425 -> 004895776 | Fernando | -5214.48
#   Row 426 -> 005324840 | Pedro    | -5233.85
#
# Delete from the bottom up so row indices for rows not yet deleted
# remain stable while we work.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(426).Delete()
$ws.Rows.Item(425).Delete()
$ws.Rows.Item(4).Delete()
